$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two test-data title strings in row 7 (C7, E7)
$ws.Range("C7").Value = "มหาศึกคนชนคน เล่ม 1"
$ws.Range("E7").Value = "มหาศึกคนชนคน เล่ม 2"

# Update row 12 height from 21 to 19.5 points
$ws.Rows.Item(12).RowHeight = 19.5
